# Auto-generated edit script for natmi Col1a2-Itgb3 LR-pairs sheet
# Adds M2 as a target cluster, making the Sending x Target matrix a full 4x4 (16 rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,7).Value = 6.423576
$ws.Cells.Item(2,8).Value = 19.270728
$ws.Cells.Item(2,9).Value = 0.001681024218962088
$ws.Cells.Item(2,10).Value = 0.001681024218962088
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,13).Value = 2.481489333333333
$ws.Cells.Item(2,14).Value = 7.444467999999999
$ws.Cells.Item(2,15).Value = 0.2345069082418988
$ws.Cells.Item(2,16).Value = 0.2345069082418987
$ws.Cells.Item(2,17).Value = 15.940035325856
$ws.Cells.Item(2,18).Value = 143.460317932704
$ws.Cells.Item(2,19).Value = 0.000394211792268552
$ws.Cells.Item(2,20).Value = 0.0003942117922685519

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,7).Value = 6.423576
$ws.Cells.Item(3,8).Value = 19.270728
$ws.Cells.Item(3,9).Value = 0.001681024218962088
$ws.Cells.Item(3,10).Value = 0.001681024218962088
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,13).Value = 7.245227
$ws.Cells.Item(3,14).Value = 21.735681
$ws.Cells.Item(3,15).Value = 0.6846919551326144
$ws.Cells.Item(3,16).Value = 0.6846919551326142
$ws.Cells.Item(3,17).Value = 46.540266271752
$ws.Cells.Item(3,18).Value = 418.862396445768
$ws.Cells.Item(3,19).Value = 0.001150983759106428
$ws.Cells.Item(3,20).Value = 0.001150983759106428

# Row 4: ECs -> M2
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,7).Value = 6.423576
$ws.Cells.Item(4,8).Value = 19.270728
$ws.Cells.Item(4,9).Value = 0.001681024218962088
$ws.Cells.Item(4,10).Value = 0.001681024218962088
$ws.Cells.Item(4,11).Value = 2.0
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.2001876666666667
$ws.Cells.Item(4,14).Value = 0.600563
$ws.Cells.Item(4,15).Value = 0.01891823194544989
$ws.Cells.Item(4,16).Value = 0.01891823194544989
$ws.Cells.Item(4,17).Value = 1.285920691096
$ws.Cells.Item(4,18).Value = 11.573286219864
$ws.Cells.Item(4,19).Value = 0.00003180200608024353
$ws.Cells.Item(4,20).Value = 0.00003180200608024353

# Row 5: ECs -> sCs
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,7).Value = 6.423576
$ws.Cells.Item(5,8).Value = 19.270728
$ws.Cells.Item(5,9).Value = 0.001681024218962088
$ws.Cells.Item(5,10).Value = 0.001681024218962088
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,13).Value = 0.6548283333333333
$ws.Cells.Item(5,14).Value = 1.964485
$ws.Cells.Item(5,15).Value = 0.06188290468003712
$ws.Cells.Item(5,16).Value = 0.06188290468003711
$ws.Cells.Item(5,17).Value = 4.20633956612
$ws.Cells.Item(5,18).Value = 37.85705609508
$ws.Cells.Item(5,19).Value = 0.0001040266615068647
$ws.Cells.Item(5,20).Value = 0.0001040266615068647

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,7).Value = 3580.644531333333
$ws.Cells.Item(6,8).Value = 10741.933594
$ws.Cells.Item(6,9).Value = 0.9370403925578976
$ws.Cells.Item(6,10).Value = 0.9370403925578976
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,13).Value = 2.481489333333333
$ws.Cells.Item(6,14).Value = 7.444467999999999
$ws.Cells.Item(6,15).Value = 0.2345069082418988
$ws.Cells.Item(6,16).Value = 0.2345069082418987
$ws.Cells.Item(6,17).Value = 8885.331210961996
$ws.Cells.Item(6,18).Value = 79967.98089865797
$ws.Cells.Item(6,19).Value = 0.2197424453565277
$ws.Cells.Item(6,20).Value = 0.2197424453565276

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,7).Value = 3580.644531333333
$ws.Cells.Item(7,8).Value = 10741.933594
$ws.Cells.Item(7,9).Value = 0.9370403925578976
$ws.Cells.Item(7,10).Value = 0.9370403925578976
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,13).Value = 7.245227
$ws.Cells.Item(7,14).Value = 21.735681
$ws.Cells.Item(7,15).Value = 0.6846919551326144
$ws.Cells.Item(7,16).Value = 0.6846919551326142
$ws.Cells.Item(7,17).Value = 25942.58243581861
$ws.Cells.Item(7,18).Value = 233483.2419223675
$ws.Cells.Item(7,19).Value = 0.6415840184186994
$ws.Cells.Item(7,20).Value = 0.6415840184186992

# Row 8: FAPs -> M2
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,7).Value = 3580.644531333333
$ws.Cells.Item(8,8).Value = 10741.933594
$ws.Cells.Item(8,9).Value = 0.9370403925578976
$ws.Cells.Item(8,10).Value = 0.9370403925578976
$ws.Cells.Item(8,11).Value = 2.0
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.2001876666666667
$ws.Cells.Item(8,14).Value = 0.600563
$ws.Cells.Item(8,15).Value = 0.01891823194544989
$ws.Cells.Item(8,16).Value = 0.01891823194544989
$ws.Cells.Item(8,17).Value = 716.8008738903802
$ws.Cells.Item(8,18).Value = 6451.207865013422
$ws.Cells.Item(8,19).Value = 0.01772714748866573
$ws.Cells.Item(8,20).Value = 0.01772714748866572

# Row 9: FAPs -> sCs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,7).Value = 3580.644531333333
$ws.Cells.Item(9,8).Value = 10741.933594
$ws.Cells.Item(9,9).Value = 0.9370403925578976
$ws.Cells.Item(9,10).Value = 0.9370403925578976
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,13).Value = 0.6548283333333333
$ws.Cells.Item(9,14).Value = 1.964485
$ws.Cells.Item(9,15).Value = 0.06188290468003712
$ws.Cells.Item(9,16).Value = 0.06188290468003711
$ws.Cells.Item(9,17).Value = 2344.707490712121
$ws.Cells.Item(9,18).Value = 21102.36741640909
$ws.Cells.Item(9,19).Value = 0.05798678129400494
$ws.Cells.Item(9,20).Value = 0.05798678129400493

# Row 10: M2 -> ECs
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,7).Value = 0.9157713333333334
$ws.Cells.Item(10,8).Value = 2.747314
$ws.Cells.Item(10,9).Value = 0.0002396537054071653
$ws.Cells.Item(10,10).Value = 0.0002396537054071653
$ws.Cells.Item(10,11).Value = 3.0
$ws.Cells.Item(10,13).Value = 2.481489333333333
$ws.Cells.Item(10,14).Value = 7.444467999999999
$ws.Cells.Item(10,15).Value = 0.2345069082418988
$ws.Cells.Item(10,16).Value = 0.2345069082418987
$ws.Cells.Item(10,17).Value = 2.272476795439111
$ws.Cells.Item(10,18).Value = 20.452291158952
$ws.Cells.Item(10,19).Value = 0.00005620044950374914
$ws.Cells.Item(10,20).Value = 0.00005620044950374913

# Row 11: M2 -> FAPs
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,7).Value = 0.9157713333333334
$ws.Cells.Item(11,8).Value = 2.747314
$ws.Cells.Item(11,9).Value = 0.0002396537054071653
$ws.Cells.Item(11,10).Value = 0.0002396537054071653
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,13).Value = 7.245227
$ws.Cells.Item(11,14).Value = 21.735681
$ws.Cells.Item(11,15).Value = 0.6846919551326144
$ws.Cells.Item(11,16).Value = 0.6846919551326142
$ws.Cells.Item(11,17).Value = 6.634971190092667
$ws.Cells.Item(11,18).Value = 59.71474071083401
$ws.Cells.Item(11,19).Value = 0.0001640889641100076
$ws.Cells.Item(11,20).Value = 0.0001640889641100075

# Row 12: M2 -> M2
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,7).Value = 0.9157713333333334
$ws.Cells.Item(12,8).Value = 2.747314
$ws.Cells.Item(12,9).Value = 0.0002396537054071653
$ws.Cells.Item(12,10).Value = 0.0002396537054071653
$ws.Cells.Item(12,11).Value = 2.0
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 0.2001876666666667
$ws.Cells.Item(12,14).Value = 0.600563
$ws.Cells.Item(12,15).Value = 0.01891823194544989
$ws.Cells.Item(12,16).Value = 0.01891823194544989
$ws.Cells.Item(12,17).Value = 0.1833261264202222
$ws.Cells.Item(12,18).Value = 1.649935137782
$ws.Cells.Item(12,19).Value = 0.000004533824385479271
$ws.Cells.Item(12,20).Value = 0.000004533824385479271

# Row 13: M2 -> sCs
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,7).Value = 0.9157713333333334
$ws.Cells.Item(13,8).Value = 2.747314
$ws.Cells.Item(13,9).Value = 0.0002396537054071653
$ws.Cells.Item(13,10).Value = 0.0002396537054071653
$ws.Cells.Item(13,11).Value = 3.0
$ws.Cells.Item(13,13).Value = 0.6548283333333333
$ws.Cells.Item(13,14).Value = 1.964485
$ws.Cells.Item(13,15).Value = 0.06188290468003712
$ws.Cells.Item(13,16).Value = 0.06188290468003711
$ws.Cells.Item(13,17).Value = 0.5996730159211111
$ws.Cells.Item(13,18).Value = 5.397057143290001
$ws.Cells.Item(13,19).Value = 0.00001483046740792931
$ws.Cells.Item(13,20).Value = 0.0000148304674079293

# Row 14: sCs -> ECs
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Col1a2"
$ws.Cells.Item(14,3).Value = "Itgb3"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3.0
$ws.Cells.Item(14,6).Value = 1.0
$ws.Cells.Item(14,7).Value = 233.243637
$ws.Cells.Item(14,8).Value = 699.7309110000001
$ws.Cells.Item(14,9).Value = 0.0610389295177331
$ws.Cells.Item(14,10).Value = 0.06103892951773311
$ws.Cells.Item(14,11).Value = 3.0
$ws.Cells.Item(14,12).Value = 1.0
$ws.Cells.Item(14,13).Value = 2.481489333333333
$ws.Cells.Item(14,14).Value = 7.444467999999999
$ws.Cells.Item(14,15).Value = 0.2345069082418988
$ws.Cells.Item(14,16).Value = 0.2345069082418987
$ws.Cells.Item(14,17).Value = 578.791597283372
$ws.Cells.Item(14,18).Value = 5209.124375550348
$ws.Cells.Item(14,19).Value = 0.01431405064359876
$ws.Cells.Item(14,20).Value = 0.01431405064359876

# Row 15: sCs -> FAPs
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Col1a2"
$ws.Cells.Item(15,3).Value = "Itgb3"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3.0
$ws.Cells.Item(15,6).Value = 1.0
$ws.Cells.Item(15,7).Value = 233.243637
$ws.Cells.Item(15,8).Value = 699.7309110000001
$ws.Cells.Item(15,9).Value = 0.0610389295177331
$ws.Cells.Item(15,10).Value = 0.06103892951773311
$ws.Cells.Item(15,11).Value = 3.0
$ws.Cells.Item(15,12).Value = 1.0
$ws.Cells.Item(15,13).Value = 7.245227
$ws.Cells.Item(15,14).Value = 21.735681
$ws.Cells.Item(15,15).Value = 0.6846919551326144
$ws.Cells.Item(15,16).Value = 0.6846919551326142
$ws.Cells.Item(15,17).Value = 1689.903096370599
$ws.Cells.Item(15,18).Value = 15209.12786733539
$ws.Cells.Item(15,19).Value = 0.04179286399069852
$ws.Cells.Item(15,20).Value = 0.04179286399069852

# Row 16: sCs -> M2
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Col1a2"
$ws.Cells.Item(16,3).Value = "Itgb3"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3.0
$ws.Cells.Item(16,6).Value = 1.0
$ws.Cells.Item(16,7).Value = 233.243637
$ws.Cells.Item(16,8).Value = 699.7309110000001
$ws.Cells.Item(16,9).Value = 0.0610389295177331
$ws.Cells.Item(16,10).Value = 0.06103892951773311
$ws.Cells.Item(16,11).Value = 2.0
$ws.Cells.Item(16,12).Value = 0.6666666666666666
$ws.Cells.Item(16,13).Value = 0.2001876666666667
$ws.Cells.Item(16,14).Value = 0.600563
$ws.Cells.Item(16,15).Value = 0.01891823194544989
$ws.Cells.Item(16,16).Value = 0.01891823194544989
$ws.Cells.Item(16,17).Value = 46.69249945587701
$ws.Cells.Item(16,18).Value = 420.232495102893
$ws.Cells.Item(16,19).Value = 0.001154748626318443
$ws.Cells.Item(16,20).Value = 0.001154748626318443

# Row 17: sCs -> sCs
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Col1a2"
$ws.Cells.Item(17,3).Value = "Itgb3"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3.0
$ws.Cells.Item(17,6).Value = 1.0
$ws.Cells.Item(17,7).Value = 233.243637
$ws.Cells.Item(17,8).Value = 699.7309110000001
$ws.Cells.Item(17,9).Value = 0.0610389295177331
$ws.Cells.Item(17,10).Value = 0.06103892951773311
$ws.Cells.Item(17,11).Value = 3.0
$ws.Cells.Item(17,12).Value = 1.0
$ws.Cells.Item(17,13).Value = 0.6548283333333333
$ws.Cells.Item(17,14).Value = 1.964485
$ws.Cells.Item(17,15).Value = 0.06188290468003712
$ws.Cells.Item(17,16).Value = 0.06188290468003711
$ws.Cells.Item(17,17).Value = 152.734542077315
$ws.Cells.Item(17,18).Value = 1374.610878695835
$ws.Cells.Item(17,19).Value = 0.003777266257117381
$ws.Cells.Item(17,20).Value = 0.003777266257117381
